$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9187528026723889
$ws.Range("B3").Value = 1.886536026929251
$ws.Range("B4").Value = 0.2058372796573881
$ws.Range("H5").Value = 97.71058970844493
$ws.Range("H6").Value = 93.54218735965073
$ws.Range("H7").Value = 96.26978552796047
$ws.Range("C8").Value = 0.1709556535546991
$ws.Range("C9").Value = 0.3667275675725769
$ws.Range("C10").Value = 0.3353156791339543
$ws.Range("D11").Value = -0.6331315114025067
$ws.Range("D12").Value = 0.2180822195408832
$ws.Range("D13").Value = -0.03931437522540901
$ws.Range("E14").Value = -0.2423051239158465
$ws.Range("E15").Value = 0.7476027935334263
$ws.Range("E16").Value = -0.7428876342878813
$ws.Range("F17").Value = -1.718306364620022
$ws.Range("F18").Value = 28.3628763538837
$ws.Range("F19").Value = -39.16836989382747
$ws.Range("G20").Value = 97.33718409582383
$ws.Range("G21").Value = 96.37610055138501
$ws.Range("G22").Value = 98.09989924092648
$ws.Range("B23").Value = -0.6190711210521657
$ws.Range("B24").Value = 0.6841272603041058
$ws.Range("H25").Value = 97.5104499833141
$ws.Range("H26").Value = 96.5701498421321
$ws.Range("C27").Value = -0.2508026900248366
$ws.Range("C28").Value = -0.1531654487793082
$ws.Range("D29").Value = 0.1298021208374179
$ws.Range("D30").Value = -0.2688808903732567
$ws.Range("E31").Value = -0.3808122002884384
$ws.Range("E32").Value = 0.8061354655834801
$ws.Range("F33").Value = 4.260618407249224
$ws.Range("F34").Value = 28.13318486538547
$ws.Range("G35").Value = 98.44083755300801
$ws.Range("G36").Value = 98.61522676209977
$ws.Range("B37").Value = -1.242244261924974
$ws.Range("B38").Value = -0.6600914489611678
$ws.Range("H39").Value = 96.08358411283866
$ws.Range("H40").Value = 90.99781590956033
$ws.Range("C41").Value = 0.7658504971464016
$ws.Range("C42").Value = 0.0391701932455986
$ws.Range("D43").Value = -0.04861296320771899
$ws.Range("D44").Value = 0.0497594565789778
$ws.Range("E45").Value = -0.08599285014618893
$ws.Range("E46").Value = 0.1243771677100176
$ws.Range("F47").Value = -0.7516248542411386
$ws.Range("F48").Value = 29.48661215871149
$ws.Range("G49").Value = 98.10516031525961
$ws.Range("G50").Value = 98.06278609418102
